$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 112 (ALC)
$ws.Range("H112").Value = 6443.875
$ws.Range("J112").Value = 7572.9
$ws.Range("L112").Value = 22718.7
$ws.Range("N112").Value = -24934.7

# Row 132 (ALC)
$ws.Range("H132").Value = 7252232.5
$ws.Range("I132").Value = 8337960.5
$ws.Range("J132").Value = 14047.333
$ws.Range("K132").Value = 25013881.5
$ws.Range("L132").Value = 42141.999
$ws.Range("M132").Value = -25011351.5
$ws.Range("N132").Value = -47201.999

# Row 137 (ALC)
$ws.Range("H137").Value = 1104.8429
$ws.Range("I137").Value = 852.383
$ws.Range("J137").Value = 1620.7391
$ws.Range("K137").Value = 2557.149
$ws.Range("L137").Value = 4862.2173
$ws.Range("M137").Value = -7.149000000000342
$ws.Range("N137").Value = -9962.2173

# Row 138 (ALC)
$ws.Range("H138").Value = 1064.0422
$ws.Range("J138").Value = 1928.2174
$ws.Range("L138").Value = 5784.6522
$ws.Range("N138").Value = -16064.6522

# Row 141 (ALC)
$ws.Range("H141").Value = 656.3570999999999
$ws.Range("I141").Value = 637.2308
$ws.Range("J141").Value = 905
$ws.Range("K141").Value = 1911.6924
$ws.Range("L141").Value = 2715
$ws.Range("M141").Value = 3268.3076
$ws.Range("N141").Value = -13075

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Range("H32").Value = 5138.467
$ws.Range("I32").Value = 4542.7734
$ws.Range("J32").Value = 9648.714
$ws.Range("K32").Value = 4542.7734
$ws.Range("L32").Value = 9648.714
$ws.Range("M32").Value = -4255.7734
$ws.Range("N32").Value = -10222.714

# Row 61 (ARM)
$ws.Range("H61").Value = 29413070
$ws.Range("I61").Value = 37038190
$ws.Range("J61").Value = 1878.2858
$ws.Range("K61").Value = 37038190
$ws.Range("L61").Value = 1878.2858
$ws.Range("M61").Value = -37037978
$ws.Range("N61").Value = -2302.2858

# Row 74 (ARM)
$ws.Range("H74").Value = 1528.4584
$ws.Range("I74").Value = 1153.6316
$ws.Range("J74").Value = 2952.8
$ws.Range("K74").Value = 1153.6316
$ws.Range("L74").Value = 2952.8
$ws.Range("M74").Value = -279.6315999999999
$ws.Range("N74").Value = -4700.8

# Row 77 (ARM)
$ws.Range("H77").Value = 1528.4584
$ws.Range("I77").Value = 1153.6316
$ws.Range("J77").Value = 2952.8
$ws.Range("K77").Value = 5768.157999999999
$ws.Range("L77").Value = 14764
$ws.Range("M77").Value = -1400.157999999999
$ws.Range("N77").Value = -23500

# Row 133 (ARM)
$ws.Range("H133").Value = 30420
$ws.Range("J133").Value = 30420
$ws.Range("L133").Value = 30420
$ws.Range("N133").Value = -35480

# Row 136 (ARM)
$ws.Range("H136").Value = 29413070
$ws.Range("I136").Value = 37038190
$ws.Range("J136").Value = 1878.2858
$ws.Range("K136").Value = 111114570
$ws.Range("L136").Value = 5634.857400000001
$ws.Range("M136").Value = -111112020
$ws.Range("N136").Value = -10734.8574

$ws = $wb.Worksheets.Item("BSM")
# Row 107 (BSM)
$ws.Range("H107").Value = 1190.9524
$ws.Range("I107").Value = 961.46155
$ws.Range("J107").Value = 1563.875
$ws.Range("K107").Value = 961.46155
$ws.Range("L107").Value = 1563.875
$ws.Range("M107").Value = 958.53845
$ws.Range("N107").Value = -5403.875

# Row 134 (BSM)
$ws.Range("H134").Value = 4312.0884
$ws.Range("I134").Value = 1084.2258
$ws.Range("J134").Value = 37666.668
$ws.Range("K134").Value = 3252.6774
$ws.Range("L134").Value = 113000.004
$ws.Range("M134").Value = -717.6773999999996
$ws.Range("N134").Value = -118070.004

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 1732.1177
$ws.Range("I31").Value = 1679.7333
$ws.Range("K31").Value = 1679.7333
$ws.Range("M31").Value = -1384.7333

# Row 34 (CRP)
$ws.Range("H34").Value = 1732.1177
$ws.Range("I34").Value = 1679.7333
$ws.Range("K34").Value = 1679.7333
$ws.Range("M34").Value = -1477.7333

# Row 58 (CRP)
$ws.Range("H58").Value = 785.5217
$ws.Range("I58").Value = 748.38464
$ws.Range("J58").Value = 992.4286
$ws.Range("K58").Value = 748.38464
$ws.Range("L58").Value = 992.4286
$ws.Range("M58").Value = -545.38464
$ws.Range("N58").Value = -1398.4286

# Row 132 (CRP)
$ws.Range("H132").Value = 4161.9756
$ws.Range("I132").Value = 4290.436
$ws.Range("J132").Value = 1657
$ws.Range("K132").Value = 12871.308
$ws.Range("L132").Value = 4971
$ws.Range("M132").Value = -10341.308
$ws.Range("N132").Value = -10031

# Row 134 (CRP)
$ws.Range("H134").Value = 13890343
$ws.Range("I134").Value = 1534.6207
$ws.Range("K134").Value = 4603.8621
$ws.Range("M134").Value = -2068.8621

# Row 135 (CRP)
$ws.Range("H135").Value = 31171.2
$ws.Range("J135").Value = 33523.555
$ws.Range("L135").Value = 33523.555
$ws.Range("N135").Value = -43663.555

# Row 136 (CRP)
$ws.Range("H136").Value = 785.5217
$ws.Range("I136").Value = 748.38464
$ws.Range("J136").Value = 992.4286
$ws.Range("K136").Value = 2245.15392
$ws.Range("L136").Value = 2977.2858
$ws.Range("M136").Value = 304.8460800000003
$ws.Range("N136").Value = -8077.2858

$ws = $wb.Worksheets.Item("CUL")
# Row 51 (CUL)
$ws.Range("H51").Value = 1852
$ws.Range("I51").Value = 1852
$ws.Range("K51").Value = 5556
$ws.Range("M51").Value = -5096

# Row 97 (CUL)
$ws.Range("H97").Value = 561.53845
$ws.Range("I97").Value = 413.75
$ws.Range("J97").Value = 798
$ws.Range("K97").Value = 1241.25
$ws.Range("L97").Value = 2394
$ws.Range("M97").Value = -745.25
$ws.Range("N97").Value = -3386

# Row 107 (CUL)
$ws.Range("H107").Value = 5626.737
$ws.Range("J107").Value = 8054.5386
$ws.Range("L107").Value = 24163.6158
$ws.Range("N107").Value = -28003.6158

$ws = $wb.Worksheets.Item("GSM")
# Row 107 (GSM)
$ws.Range("H107").Value = 540.93335
$ws.Range("I107").Value = 743.4
$ws.Range("J107").Value = 338.46667
$ws.Range("K107").Value = 743.4
$ws.Range("L107").Value = 338.46667
$ws.Range("M107").Value = 1176.6
$ws.Range("N107").Value = -4178.46667

# Row 126 (GSM)
$ws.Range("H126").Value = 2831.3076
$ws.Range("I126").Value = 1850
$ws.Range("J126").Value = 3672.4285
$ws.Range("K126").Value = 5550
$ws.Range("L126").Value = 11017.2855
$ws.Range("M126").Value = -3080
$ws.Range("N126").Value = -15957.2855

# Row 134 (GSM)
$ws.Range("H134").Value = 28521.428
$ws.Range("J134").Value = 28521.428
$ws.Range("L134").Value = 85564.284
$ws.Range("N134").Value = -90634.284

$ws = $wb.Worksheets.Item("LTW")
# Row 45 (LTW)
$ws.Range("H45").Value = 5000
$ws.Range("I45").Value = 5000
$ws.Range("K45").Value = 5000
$ws.Range("M45").Value = -4593

# Row 61 (LTW)
$ws.Range("H61").Value = 1162.9
$ws.Range("I61").Value = 1024.8889
$ws.Range("K61").Value = 1024.8889
$ws.Range("M61").Value = -822.8888999999999

# Row 113 (LTW)
$ws.Range("H113").Value = 1162.9
$ws.Range("I113").Value = 1024.8889
$ws.Range("K113").Value = 1024.8889
$ws.Range("M113").Value = 1145.1111

# Row 132 (LTW)
$ws.Range("H132").Value = 21609.18
$ws.Range("I132").Value = 1468.7142
$ws.Range("K132").Value = 4406.142599999999
$ws.Range("M132").Value = -1876.142599999999

# Row 136 (LTW)
$ws.Range("H136").Value = 3367.4348
$ws.Range("I136").Value = 3224.932
$ws.Range("J136").Value = 6502.5
$ws.Range("K136").Value = 9674.795999999998
$ws.Range("L136").Value = 19507.5
$ws.Range("M136").Value = -7124.795999999998
$ws.Range("N136").Value = -24607.5

$ws = $wb.Worksheets.Item("WVR")
# Row 6 (WVR)
$ws.Range("H6").Value = 720
$ws.Range("I6").Value = 130
$ws.Range("J6").Value = 916.6667
$ws.Range("K6").Value = 130
$ws.Range("L6").Value = 916.6667
$ws.Range("M6").Value = -15
$ws.Range("N6").Value = -1146.6667

# Row 48 (WVR)
$ws.Range("H48").Value = 7516.25
$ws.Range("J48").Value = 7516.25
$ws.Range("L48").Value = 7516.25
$ws.Range("N48").Value = -8654.25

# Row 132 (WVR)
$ws.Range("H132").Value = 4263.5
$ws.Range("I132").Value = 6318.3335
$ws.Range("J132").Value = 1181.25
$ws.Range("K132").Value = 18955.0005
$ws.Range("L132").Value = 3543.75
$ws.Range("M132").Value = -16425.0005
$ws.Range("N132").Value = -8603.75

